$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the order dates in rows 12-17: they all become the same date (40940) ---
$ws.Cells.Item(12,1).Value = 40940
$ws.Cells.Item(13,1).Value = 40940
$ws.Cells.Item(14,1).Value = 40940
$ws.Cells.Item(15,1).Value = 40940
$ws.Cells.Item(16,1).Value = 40940
$ws.Cells.Item(17,1).Value = 40940

# --- Rows 18-26: dates move from 40947 to 40942 ---
$ws.Cells.Item(18,1).Value = 40942
$ws.Cells.Item(19,1).Value = 40942
$ws.Cells.Item(20,1).Value = 40942
$ws.Cells.Item(21,1).Value = 40942
$ws.Cells.Item(22,1).Value = 40942
$ws.Cells.Item(23,1).Value = 40942
$ws.Cells.Item(24,1).Value = 40942
$ws.Cells.Item(25,1).Value = 40942
$ws.Cells.Item(26,1).Value = 40942

# B20-B26 keep their original style, just update the value
$ws.Cells.Item(20,2).Value = 40942
$ws.Cells.Item(21,2).Value = 40942
$ws.Cells.Item(22,2).Value = 40942
$ws.Cells.Item(23,2).Value = 40942
$ws.Cells.Item(24,2).Value = 40942
$ws.Cells.Item(25,2).Value = 40942
$ws.Cells.Item(26,2).Value = 40942

# B18 and B19 get re-typed with an explicit black font colour, which Excel
# records as a brand-new cell style (keeping the existing short-date format)
$ws.Cells.Item(18,2).Style = "Normal"
$ws.Cells.Item(18,2).Value = 40942
$ws.Cells.Item(18,2).Font.Color = 0
$ws.Cells.Item(18,2).NumberFormat = "mm-dd-yy"

# Re-use the style just created for B18 on B19 as well
$ws.Cells.Item(18,2).Copy() | Out-Null
$ws.Cells.Item(19,2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(19,2).Value = 40942

# --- Row 26: the F column used to (mis-)point at the stray "v" string;
#     point it at "Tornicalvo" like all the other screw rows ---
$ws.Cells.Item(26,6).Value = "Tornicalvo"

# --- Add the missing row 27 for the HeatBed (this reclaims the now orphaned
#     "v" shared-string slot and turns it into "HeatBed") ---
$ws.Cells.Item(27,1).Value = 40943
$ws.Cells.Item(26,1).Copy() | Out-Null
$ws.Cells.Item(27,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(27,1).Value = 40943
$ws.Cells.Item(27,3).Value = 1
$ws.Cells.Item(27,4).Value = "HeatBed"
$ws.Cells.Item(27,5).Value = 36.5
$ws.Cells.Item(27,6).Value = "http://reprapworld.com/?products_details&products_id=121"

$ws.Cells.Item(18,2).Select() | Out-Null
